$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("order")

# The order sheet is being hooked up to read from a CSV that also carries a
# Postal Code field. Make room for it as a new column between Province and
# Credit Card Number by shifting Credit Card Number / CCV / Pizza Type one
# column to the right (values + formatting), then fill in the new column.
$ws.Range("I1:I2").Copy($ws.Range("J1:J2"))
$ws.Range("H1:H2").Copy($ws.Range("I1:I2"))
$ws.Range("G1:G2").Copy($ws.Range("H1:H2"))

$ws.Range("G1").Value = "Postal Code"
$ws.Range("G2").Value = "L6A 1H6"

# Re-fit the (now-shifted) columns so their widths track their new content,
# matching how Excel auto-recalculates "best fit" column widths.
$ws.Columns("G").ColumnWidth = 11.42578125
$ws.Columns("H").ColumnWidth = 18.85546875
$ws.Columns("I").ColumnWidth = 4.5703125
$ws.Columns("J").ColumnWidth = 10.28515625

# Update the active selection to match where the user was working.
$ws.Range("F6").Select()
